# Generate Report for Handoff
# Updates the "Ready for handoff" rows (4-7) on both the zh-cn and de-de
# localization-status sheets: the Priority moves from "low" to "ht" and the
# Latest Handoff Datetime is refreshed to a newer timestamp (30 seconds
# later than the previous value) to reflect the freshly generated handoff
# report.

$wb = $excel.ActiveWorkbook

$ws_overview = $wb.Worksheets.Item("Overview")
$ws_zh = $wb.Worksheets.Item("zh-cn")
$ws_de = $wb.Worksheets.Item("de-de")

for ($row = 4; $row -le 7; $row++) {
    $ws_zh.Cells.Item($row, 5).Value = "ht"
    $ws_de.Cells.Item($row, 5).Value = "ht"
}

for ($row = 4; $row -le 7; $row++) {
    $cell = $ws_zh.Range("H$row")
    $cell.Value = "2016-08-13 14:36:14"
    $cell.NumberFormat = "yyyy-mm-dd HH:mm:ss"

    $cell = $ws_de.Range("H$row")
    $cell.Value = "2016-08-13 14:36:22"
    $cell.NumberFormat = "yyyy-mm-dd HH:mm:ss"
}

# The Overview sheet's "Latest HO Xliff Generate Date" tracks the newest
# handoff datetime across languages, which is now the de-de timestamp.
for ($row = 4; $row -le 7; $row++) {
    $cell = $ws_overview.Range("G$row")
    $cell.Value = "2016-08-13 14:36:22"
    $cell.NumberFormat = "yyyy-mm-dd HH:mm:ss"
}
